# Updates cryptocurrency price/volume data to match latest scrape
# All D/E columns are text-formatted (prices use "." as thousands
# separators in some rows, so the column is text, not numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their existing Text number format so that
# values such as "214.99" or "1.00" are stored as strings, matching the
# original workbook contents (inline strings), not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.839.08"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.633.57"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "214.99"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "19.94"
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.635.02"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.859.30"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "0.560"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "63.05"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "25.845.56"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "193.75"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "1.75"
$ws.Range("E25").Value = "  -5.04%  "
$ws.Range("D26").Value = "138.37"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("E27").Value = "  -4.67%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.120.31"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.547"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "0.0156"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "99.35"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").Value = "0.799"
$ws.Range("D45").Value = "0.0₆0110"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("E47").Value = "  -5.03%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.57"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("E51").Value = "  -0.26%  "
